$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 119. This shifts the existing rows 119:196 down to
# 120:197, preserving all their data and formatting (this is how the source
# dataset records a newly-added weekly price entry at the top of the
# historical list for this market/product).
$ws.Rows("119:119").Insert()

# Populate the newly inserted row 119 with the new weekly record.
$ws.Cells.Item(119, 1).Value = 11
$ws.Cells.Item(119, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(119, 3).Value = "Bíobío"
$ws.Cells.Item(119, 4).Value = 44572
$ws.Cells.Item(119, 5).Value = 8
$ws.Cells.Item(119, 6).Value = 100112008
$ws.Cells.Item(119, 7).Value = "Coliflor"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 220
$ws.Cells.Item(119, 11).Value = 700
$ws.Cells.Item(119, 12).Value = 800
$ws.Cells.Item(119, 13).Value = 755
$ws.Cells.Item(119, 14).Value = "`$/unidad"
$ws.Cells.Item(119, 15).Value = "Región Metropolitana"
$ws.Cells.Item(119, 16).Value = 755
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"
